# Append two new案件 (job listing) rows scraped at 2025-11-09 01:23:23,
# and refresh the timestamp / top-row figures of the existing listings
# on the "ランサーズ" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-09 01:23:23"

# --- Refresh the "取得日時" (fetched-at) column for every existing row ---
$ws.Range("A2:A11").Value = $newTimestamp

# --- Row 2's price range / priority score moved between scrapes ---
$ws.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("G2").Value = 458

# --- Column D ("価格") got a bit narrower in this edit ---
$ws.Columns.Item(4).ColumnWidth = 27.15

# --- Append the two newly discovered listings as rows 12 and 13 ---
$ws.Range("A12").Value = $newTimestamp
$ws.Range("B12").Value = "MT4 RSXを使用したEAの作成依頼"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5430008"
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5430008")
$ws.Range("F12").Style = "Hyperlink"
$ws.Range("G12").Value = 10

$ws.Range("A13").Value = $newTimestamp
$ws.Range("B13").Value = "【急募】LINE × QRコード連携で自動取得設定を実現!"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5430015"
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5430015")
$ws.Range("F13").Style = "Hyperlink"
$ws.Range("G13").Value = 10
